# Updates cryptos list figures (price / 1h volume change) and fixes a
# swapped row (ApeXProtocol <-> ThetaToken), matching the upstream GitHub
# Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" values look like plain decimal numbers (e.g.
# "605.65"), which Excel would otherwise auto-convert to a numeric cell.
# The source data keeps these as literal text (inline/shared strings), so
# we temporarily force a Text number format on those specific cells,
# assign the string value, then restore the default "Normal" style so the
# cell's style index goes back to the sheet's default (no visible
# formatting change is introduced).
$textForceCells = @("D5", "D6", "D7", "D11", "D13", "D15", "D16", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D44", "D45", "D46", "D48", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price column updates
$ws.Range("D2").Value  = "70.224.05"
$ws.Range("D3").Value  = "3.622.12"
$ws.Range("D5").Value  = "605.65"
$ws.Range("D6").Value  = "196.78"
$ws.Range("D7").Value  = "0.628"
$ws.Range("D11").Value = "54.24"
$ws.Range("D13").Value = "9.56"
$ws.Range("D14").Value = "4.201.64"
$ws.Range("D15").Value = "13.25"
$ws.Range("D16").Value = "594.06"
$ws.Range("D18").Value = "70.362.57"
$ws.Range("D19").Value = "3.622.99"
$ws.Range("D21").Value = "0.998"
$ws.Range("D22").Value = "17.72"
$ws.Range("D23").Value = "5.16"
$ws.Range("D24").Value = "102.75"
$ws.Range("D25").Value = "4.63"
$ws.Range("D26").Value = "3.05"
$ws.Range("D27").Value = "10.83"
$ws.Range("D28").Value = "9.61"
$ws.Range("D29").Value = "34.04"
$ws.Range("D30").Value = "4.67"
$ws.Range("D31").Value = "7.23"
$ws.Range("D33").Value = "0.118"
$ws.Range("D34").Value = "63.27"
$ws.Range("D36").Value = "3.935.87"
$ws.Range("D37").Value = "3.19"
$ws.Range("D38").Value = "523.30"
$ws.Range("D39").Value = "1.00"
$ws.Range("D40").Value = "37.40"
$ws.Range("D41").Value = "0.394"
$ws.Range("D44").Value = "0.0457"
$ws.Range("D45").Value = "2.87"
$ws.Range("D46").Value = "3.35"
$ws.Range("D48").Value = "8.64"
$ws.Range("D50").Value = "0.000254"

# Volume(1h) column updates
$ws.Range("E3").Value  = "  +3.43%  "
$ws.Range("E4").Value  = "  -0.06%  "
$ws.Range("E5").Value  = "  +1.17%  "
$ws.Range("E6").Value  = "  +0.59%  "
$ws.Range("E7").Value  = "  +0.75%  "
$ws.Range("E9").Value  = "  -0.23%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  -1.21%  "
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("E33").Value = "  +2.42%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  +11.24%  "
$ws.Range("E36").Value = "  +5.31%  "
$ws.Range("E37").Value = "  +5.20%  "
$ws.Range("E38").Value = "  +3.90%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  +1.91%  "
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("E50").Value = "  +4.93%  "
$ws.Range("E51").Value = "  +3.94%  "

# Rows 45 & 46 swapped coin identity (Coin name + Link)
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"

# Restore the default cell style on cells we temporarily forced to Text
# number format above (keeps them as text values, but without leaving a
# stray number-format/style applied to the cell).
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
